$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ur = $ws.UsedRange
$lastRow = $ur.Row + $ur.Rows.Count - 1
if ($lastRow -lt 431) { $lastRow = 431 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
